# Auto-generated edit script applying the scheduled-runner profit recalculation
# to the Lamia_Profits workbook's per-leve profit columns (H:N) across all 8 sheets.
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
# ALC row 15
$ws.Range("H15").Value = 1351.7576
$ws.Range("I15").Value = 1351.7576
$ws.Range("K15").Value = 4055.2728
$ws.Range("M15").Value = -3886.2728
# ALC row 18
$ws.Range("H18").Value = 1289.2
$ws.Range("I18").Value = 1486.5
$ws.Range("J18").Value = 500
$ws.Range("K18").Value = 1486.5
$ws.Range("L18").Value = 500
$ws.Range("M18").Value = -1202.5
$ws.Range("N18").Value = -1068
# ALC row 20
$ws.Range("H20").Value = 5000
$ws.Range("I20").Value = 5000
$ws.Range("K20").Value = 5000
$ws.Range("M20").Value = -4770
# ALC row 35
$ws.Range("H35").Value = 5000
$ws.Range("I35").Value = 5000
$ws.Range("K35").Value = 5000
$ws.Range("M35").Value = -4621
# ALC row 80
$ws.Range("H80").Value = 2535.5
$ws.Range("I80").Value = 279
$ws.Range("J80").Value = 3789.111
$ws.Range("K80").Value = 837
$ws.Range("L80").Value = 11367.333
$ws.Range("M80").Value = 161
$ws.Range("N80").Value = -13363.333
# ALC row 83
$ws.Range("H83").Value = 2535.5
$ws.Range("I83").Value = 279
$ws.Range("J83").Value = 3789.111
$ws.Range("K83").Value = 2511
$ws.Range("L83").Value = 34101.999
$ws.Range("M83").Value = 2481
$ws.Range("N83").Value = -44085.999
# ALC row 98
$ws.Range("H98").Value = 1969.125
$ws.Range("I98").Value = 1107.5714
$ws.Range("K98").Value = 1107.5714
$ws.Range("M98").Value = 390.4286
# ALC row 106
$ws.Range("H106").Value = 10272.625
$ws.Range("I106").Value = 1591.5
$ws.Range("J106").Value = 15481.3
$ws.Range("K106").Value = 1591.5
$ws.Range("L106").Value = 15481.3
$ws.Range("M106").Value = -960.5
$ws.Range("N106").Value = -16743.3
# ALC row 122
$ws.Range("H122").Value = 1969.125
$ws.Range("I122").Value = 1107.5714
$ws.Range("K122").Value = 3322.7142
$ws.Range("M122").Value = -872.7142000000003
# ALC row 132
$ws.Range("H132").Value = 1676.6774
$ws.Range("I132").Value = 1463.25
$ws.Range("K132").Value = 4389.75
$ws.Range("M132").Value = -1859.75
# ALC row 133
$ws.Range("H133").Value = 65166
$ws.Range("J133").Value = 65166
$ws.Range("L133").Value = 65166
$ws.Range("N133").Value = -75286
# ALC row 135
$ws.Range("H135").Value = 970.4375
$ws.Range("I135").Value = 823.2143
$ws.Range("K135").Value = 7408.928699999999
$ws.Range("M135").Value = -4873.928699999999
# ALC row 141
$ws.Range("H141").Value = 5332.9644
$ws.Range("I141").Value = 2387.7144
$ws.Range("K141").Value = 7163.1432
$ws.Range("M141").Value = -1983.1432

$ws = $wb.Worksheets.Item("ARM")
# ARM row 122
$ws.Range("H122").Value = 7666.6665
$ws.Range("I122").Value = 0
$ws.Range("J122").Value = 7666.6665
$ws.Range("K122").Value = 0
$ws.Range("L122").Value = 22999.9995
$ws.Range("M122").ClearContents()
$ws.Range("N122").Value = -27899.9995

$ws = $wb.Worksheets.Item("BSM")
# BSM row 107
$ws.Range("H107").Value = 629.5714
$ws.Range("I107").Value = 625.0909
$ws.Range("K107").Value = 625.0909
$ws.Range("M107").Value = 1294.9091

$ws = $wb.Worksheets.Item("CRP")
# CRP row 7
$ws.Range("H7").Value = 279.3846
$ws.Range("I7").Value = 251.375
$ws.Range("J7").Value = 324.2
$ws.Range("K7").Value = 251.375
$ws.Range("L7").Value = 324.2
$ws.Range("M7").Value = -138.375
$ws.Range("N7").Value = -550.2
# CRP row 32
$ws.Range("H32").Value = 2592.8572
$ws.Range("I32").Value = 1033.3334
$ws.Range("J32").Value = 3762.5
$ws.Range("K32").Value = 1033.3334
$ws.Range("L32").Value = 3762.5
$ws.Range("M32").Value = -717.3334
$ws.Range("N32").Value = -4394.5
# CRP row 58
$ws.Range("H58").Value = 8467.200000000001
$ws.Range("I58").Value = 6165.3335
$ws.Range("K58").Value = 6165.3335
$ws.Range("M58").Value = -5962.3335
# CRP row 122
$ws.Range("H122").Value = 7945.75
$ws.Range("I122").Value = 2642.7144
$ws.Range("J122").Value = 15370
$ws.Range("K122").Value = 7928.1432
$ws.Range("L122").Value = 46110
$ws.Range("M122").Value = -5478.1432
$ws.Range("N122").Value = -51010
# CRP row 136
$ws.Range("H136").Value = 8467.200000000001
$ws.Range("I136").Value = 6165.3335
$ws.Range("K136").Value = 18496.0005
$ws.Range("M136").Value = -15946.0005

$ws = $wb.Worksheets.Item("CUL")
# CUL row 41
$ws.Range("H41").Value = 0
$ws.Range("J41").Value = 0
$ws.Range("L41").Value = 0
$ws.Range("N41").ClearContents()
# CUL row 48
$ws.Range("H48").Value = 0
$ws.Range("I48").Value = 0
$ws.Range("K48").Value = 0
$ws.Range("M48").ClearContents()

$ws = $wb.Worksheets.Item("GSM")
# GSM row 21
$ws.Range("H21").Value = 0
$ws.Range("I21").Value = 0
$ws.Range("K21").Value = 0
$ws.Range("M21").ClearContents()
# GSM row 30
$ws.Range("H30").Value = 0
$ws.Range("I30").Value = 0
$ws.Range("K30").Value = 0
$ws.Range("M30").ClearContents()
# GSM row 80
$ws.Range("H80").Value = 4457.2
$ws.Range("I80").Value = 1295.75
$ws.Range("J80").Value = 5606.8184
$ws.Range("K80").Value = 1295.75
$ws.Range("L80").Value = 5606.8184
$ws.Range("M80").Value = -297.75
$ws.Range("N80").Value = -7602.8184
# GSM row 83
$ws.Range("H83").Value = 4457.2
$ws.Range("I83").Value = 1295.75
$ws.Range("J83").Value = 5606.8184
$ws.Range("K83").Value = 6478.75
$ws.Range("L83").Value = 28034.092
$ws.Range("M83").Value = -1486.75
$ws.Range("N83").Value = -38018.092
# GSM row 122
$ws.Range("H122").Value = 8982.799999999999
$ws.Range("I122").Value = 4953.5
$ws.Range("K122").Value = 14860.5
$ws.Range("M122").Value = -12410.5
# GSM row 138
$ws.Range("H138").Value = 50000
$ws.Range("J138").Value = 50000
$ws.Range("L138").Value = 50000
$ws.Range("N138").Value = -60280

$ws = $wb.Worksheets.Item("LTW")
# LTW row 40
$ws.Range("H40").Value = 2946
$ws.Range("I40").Value = 2304.1904
$ws.Range("K40").Value = 2304.1904
$ws.Range("M40").Value = -2168.1904
# LTW row 82
$ws.Range("H82").Value = 2983.8438
$ws.Range("I82").Value = 857.15
$ws.Range("J82").Value = 6528.3335
$ws.Range("K82").Value = 857.15
$ws.Range("L82").Value = 6528.3335
$ws.Range("M82").Value = -496.15
$ws.Range("N82").Value = -7250.3335
# LTW row 85
$ws.Range("H85").Value = 2983.8438
$ws.Range("I85").Value = 857.15
$ws.Range("J85").Value = 6528.3335
$ws.Range("K85").Value = 857.15
$ws.Range("L85").Value = 6528.3335
$ws.Range("M85").Value = 390.85
$ws.Range("N85").Value = -9024.333500000001
# LTW row 100
$ws.Range("H100").Value = 5262.2144
$ws.Range("I100").Value = 3777.3
$ws.Range("K100").Value = 3777.3
$ws.Range("M100").Value = -3236.3
# LTW row 122
$ws.Range("H122").Value = 4975.276
$ws.Range("I122").Value = 4499.16
$ws.Range("J122").Value = 7951
$ws.Range("K122").Value = 13497.48
$ws.Range("L122").Value = 23853
$ws.Range("M122").Value = -11047.48
$ws.Range("N122").Value = -28753
# LTW row 132
$ws.Range("H132").Value = 3327.0715
$ws.Range("I132").Value = 1960.8572
$ws.Range("K132").Value = 5882.571599999999
$ws.Range("M132").Value = -3352.571599999999

$ws = $wb.Worksheets.Item("WVR")
# WVR row 98
$ws.Range("H98").Value = 50000
$ws.Range("J98").Value = 50000
$ws.Range("L98").Value = 50000
$ws.Range("N98").Value = -55990
# WVR row 108
$ws.Range("H108").Value = 50000
$ws.Range("J108").Value = 50000
$ws.Range("L108").Value = 50000
$ws.Range("N108").Value = -57680
# WVR row 122
$ws.Range("H122").Value = 17584.166
$ws.Range("I122").Value = 10000
$ws.Range("K122").Value = 30000
$ws.Range("M122").Value = -27550
# WVR row 136
$ws.Range("H136").Value = 7541.143
$ws.Range("I136").Value = 3395.75
$ws.Range("K136").Value = 10187.25
$ws.Range("M136").Value = -7637.25
